# Update the "repaymentstrategy" value on the ProductLoanInput sheet
# from "RBI (India)" to "Overdue/Due Fee/Int,Principal", and leave the
# selection on that cell (B17), matching the author's edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ProductLoanInput")

$ws.Range("B17").Value = "Overdue/Due Fee/Int,Principal"

$ws.Activate()
$ws.Range("B17").Select()
